$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H105").Value = 32695.285
$ws.Range("J105").Value = 32695.285
$ws.Range("L105").Value = 32695.285
$ws.Range("N105").Value = -39683.285

$ws.Range("H109").Value = 33890
$ws.Range("J109").Value = 33890
$ws.Range("L109").Value = 33890
$ws.Range("N109").Value = -36664

$ws.Range("H130").Value = 43726.668
$ws.Range("J130").Value = 43726.668
$ws.Range("L130").Value = 43726.668
$ws.Range("N130").Value = -53766.668

$ws.Range("H132").Value = 37039388
$ws.Range("I132").Value = 43479976
$ws.Range("J132").Value = 6001.5
$ws.Range("K132").Value = 130439928
$ws.Range("L132").Value = 18004.5
$ws.Range("M132").Value = -130437398
$ws.Range("N132").Value = -23064.5

$ws.Range("H137").Value = 2270439
$ws.Range("I137").Value = 2802954
$ws.Range("J137").Value = 7250
$ws.Range("K137").Value = 8408862
$ws.Range("L137").Value = 21750
$ws.Range("M137").Value = -8406312
$ws.Range("N137").Value = -26850


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1747.5714
$ws.Range("I61").Value = 1660.6
$ws.Range("J61").Value = 1965
$ws.Range("K61").Value = 1660.6
$ws.Range("L61").Value = 1965
$ws.Range("M61").Value = -1448.6
$ws.Range("N61").Value = -2389

$ws.Range("H74").Value = 1942.3572
$ws.Range("I74").Value = 759.3
$ws.Range("J74").Value = 4900
$ws.Range("K74").Value = 759.3
$ws.Range("L74").Value = 4900
$ws.Range("M74").Value = 114.7
$ws.Range("N74").Value = -6648

$ws.Range("H77").Value = 1942.3572
$ws.Range("I77").Value = 759.3
$ws.Range("J77").Value = 4900
$ws.Range("K77").Value = 3796.5
$ws.Range("L77").Value = 24500
$ws.Range("M77").Value = 571.5
$ws.Range("N77").Value = -33236

$ws.Range("H93").Value = 24000
$ws.Range("J93").Value = 24000
$ws.Range("L93").Value = 24000
$ws.Range("N93").Value = -28992

$ws.Range("H97").Value = 1795.6154
$ws.Range("I97").Value = 1030.1818
$ws.Range("K97").Value = 1030.1818
$ws.Range("M97").Value = -534.1818000000001

$ws.Range("H123").Value = 50079.715
$ws.Range("J123").Value = 50079.715
$ws.Range("L123").Value = 50079.715
$ws.Range("N123").Value = -59879.715

$ws.Range("H132").Value = 3080.8
$ws.Range("I132").Value = 1342.8572
$ws.Range("J132").Value = 4601.5
$ws.Range("K132").Value = 4028.5716
$ws.Range("L132").Value = 13804.5
$ws.Range("M132").Value = -1498.5716
$ws.Range("N132").Value = -18864.5

$ws.Range("H136").Value = 1747.5714
$ws.Range("I136").Value = 1660.6
$ws.Range("J136").Value = 1965
$ws.Range("K136").Value = 4981.799999999999
$ws.Range("L136").Value = 5895
$ws.Range("M136").Value = -2431.799999999999
$ws.Range("N136").Value = -10995

$ws.Range("H137").Value = 41530
$ws.Range("J137").Value = 41530
$ws.Range("L137").Value = 41530
$ws.Range("N137").Value = -51730


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H10").Value = 7777
$ws.Range("J10").Value = 7777
$ws.Range("L10").Value = 7777
$ws.Range("N10").Value = -8057

$ws.Range("H11").Value = 2287.4
$ws.Range("I11").Value = 866.3333
$ws.Range("J11").Value = 2538.1765
$ws.Range("K11").Value = 866.3333
$ws.Range("L11").Value = 2538.1765
$ws.Range("M11").Value = -726.3333
$ws.Range("N11").Value = -2818.1765

$ws.Range("H24").Value = 3003.2
$ws.Range("I24").Value = 3004
$ws.Range("J24").Value = 3000
$ws.Range("K24").Value = 3004
$ws.Range("L24").Value = 3000
$ws.Range("M24").Value = -2769
$ws.Range("N24").Value = -3470

$ws.Range("H134").Value = 4120
$ws.Range("I134").Value = 1944
$ws.Range("J134").Value = 15000
$ws.Range("K134").Value = 5832
$ws.Range("L134").Value = 45000
$ws.Range("M134").Value = -3297
$ws.Range("N134").Value = -50070

$ws.Range("H137").Value = 50750
$ws.Range("J137").Value = 50750
$ws.Range("L137").Value = 50750
$ws.Range("N137").Value = -60950


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 34166.5
$ws.Range("J52").Value = 34166.5
$ws.Range("L52").Value = 34166.5
$ws.Range("N52").Value = -34754.5

$ws.Range("H120").Value = 30264.25
$ws.Range("I120").Value = 15000
$ws.Range("J120").Value = 35352.332
$ws.Range("K120").Value = 15000
$ws.Range("L120").Value = 35352.332
$ws.Range("M120").Value = -11371
$ws.Range("N120").Value = -42610.332

$ws.Range("H137").Value = 32084.75
$ws.Range("J137").Value = 32084.75
$ws.Range("L137").Value = 32084.75
$ws.Range("N137").Value = -42284.75

$ws.Range("H139").Value = 38899.668
$ws.Range("J139").Value = 38899.668
$ws.Range("L139").Value = 38899.668
$ws.Range("N139").Value = -49179.668


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 17246.545
$ws.Range("J39").Value = 18921.2
$ws.Range("L39").Value = 56763.60000000001
$ws.Range("N39").Value = -57351.60000000001

$ws.Range("H131").Value = 774.91
$ws.Range("I131").Value = 317.5
$ws.Range("J131").Value = 825.73334
$ws.Range("K131").Value = 952.5
$ws.Range("L131").Value = 2477.20002
$ws.Range("M131").Value = 4087.5
$ws.Range("N131").Value = -12557.20002


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3672.92
$ws.Range("I122").Value = 3443.3572
$ws.Range("J122").Value = 3965.0908
$ws.Range("K122").Value = 10330.0716
$ws.Range("L122").Value = 11895.2724
$ws.Range("M122").Value = -7880.071599999999
$ws.Range("N122").Value = -16795.2724

$ws.Range("H123").Value = 10322.75
$ws.Range("J123").Value = 10322.75
$ws.Range("L123").Value = 10322.75
$ws.Range("N123").Value = -15222.75

$ws.Range("H126").Value = 3388.06
$ws.Range("J126").Value = 4871.76
$ws.Range("L126").Value = 14615.28
$ws.Range("N126").Value = -19555.28

$ws.Range("H132").Value = 5552.4287
$ws.Range("I132").Value = 4612.364
$ws.Range("K132").Value = 13837.092
$ws.Range("M132").Value = -11307.092

$ws.Range("H137").Value = 40206.4
$ws.Range("J137").Value = 40206.4
$ws.Range("L137").Value = 40206.4
$ws.Range("N137").Value = -50406.4


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1921.0526
$ws.Range("I22").Value = 1681.6364
$ws.Range("J22").Value = 2250.25
$ws.Range("K22").Value = 1681.6364
$ws.Range("L22").Value = 2250.25
$ws.Range("M22").Value = -1386.6364
$ws.Range("N22").Value = -2840.25

$ws.Range("H27").Value = 1921.0526
$ws.Range("I27").Value = 1681.6364
$ws.Range("J27").Value = 2250.25
$ws.Range("K27").Value = 1681.6364
$ws.Range("L27").Value = 2250.25
$ws.Range("M27").Value = -1574.6364
$ws.Range("N27").Value = -2464.25

$ws.Range("H46").Value = 1666.7333
$ws.Range("I46").Value = 881.375
$ws.Range("J46").Value = 2564.2856
$ws.Range("K46").Value = 881.375
$ws.Range("L46").Value = 2564.2856
$ws.Range("M46").Value = -693.375
$ws.Range("N46").Value = -2940.2856

$ws.Range("H82").Value = 4618.5
$ws.Range("I82").Value = 4791
$ws.Range("J82").Value = 3756
$ws.Range("K82").Value = 4791
$ws.Range("L82").Value = 3756
$ws.Range("M82").Value = -4430
$ws.Range("N82").Value = -4478

$ws.Range("H85").Value = 4618.5
$ws.Range("I85").Value = 4791
$ws.Range("J85").Value = 3756
$ws.Range("K85").Value = 4791
$ws.Range("L85").Value = 3756
$ws.Range("M85").Value = -3543
$ws.Range("N85").Value = -6252

$ws.Range("H132").Value = 5378.909
$ws.Range("I132").Value = 2953.4666
$ws.Range("J132").Value = 7400.1113
$ws.Range("K132").Value = 8860.399800000001
$ws.Range("L132").Value = 22200.3339
$ws.Range("M132").Value = -6330.399800000001
$ws.Range("N132").Value = -27260.3339


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 39912.5
$ws.Range("J123").Value = 39912.5
$ws.Range("L123").Value = 39912.5
$ws.Range("N123").Value = -49712.5

$ws.Range("H126").Value = 3374.0588
$ws.Range("I126").Value = 2269.4285
$ws.Range("J126").Value = 4147.3
$ws.Range("K126").Value = 6808.2855
$ws.Range("L126").Value = 12441.9
$ws.Range("M126").Value = -4338.2855
$ws.Range("N126").Value = -17381.9

$ws.Range("H132").Value = 13339466
$ws.Range("I132").Value = 10136.363
$ws.Range("J132").Value = 23812512
$ws.Range("K132").Value = 30409.089
$ws.Range("L132").Value = 71437536
$ws.Range("M132").Value = -27879.089
$ws.Range("N132").Value = -71442596

$ws.Range("H136").Value = 4213.6313
$ws.Range("I136").Value = 3548.262
$ws.Range("K136").Value = 10644.786
$ws.Range("M136").Value = -8094.786

$ws.Range("H141").Value = 35665
$ws.Range("J141").Value = 35665
$ws.Range("L141").Value = 35665
$ws.Range("N141").Value = -46025

